# Remove the duplicate "AFTER PARTY / SNRS / Dortmund" entry that sits at
# row 595 (an exact duplicate of the entry that remains further down the
# list). Deleting the entire row shifts every following row up by one,
# shrinking the sheet's used range from A1:E618 to A1:E617.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(595).Delete()
